# Rotate the four reference URLs shown on the "References" slides.
# Old order: en.wikipedia.org -> www.nih.gov -> scholar.google.com -> www.jstor.org
# New order: www.nih.gov -> scholar.google.com -> www.jstor.org -> en.wikipedia.org
# (each URL moves up one slot, wikipedia wraps around to the bottom)

$p = $ppt.ActivePresentation

$newUrls = @(
    "https://www.nih.gov/",
    "https://scholar.google.com/",
    "https://www.jstor.org/",
    "https://en.wikipedia.org/wiki/Main_Page"
)

$slideIndexes = @(7, 13, 19, 25, 29, 33)

foreach ($idx in $slideIndexes) {
    $s = $p.Slides.Item($idx)
    $sh = $s.Shapes.Item(2)
    $tr = $sh.TextFrame.TextRange

    for ($i = 0; $i -lt $newUrls.Length; $i++) {
        $paraIndex = $i + 2
        $para = $tr.Paragraphs($paraIndex, 1)
        # Clear first so the new text is inserted as a single fresh run
        # instead of being merged/diffed against the old run's text.
        $para.Text = ""
        $para.Text = $newUrls[$i]
    }
}
